$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (D) and Volume(1h) (E) columns with latest
# scraped values, preserving the original text formatting (leading apostrophe
# forces Excel to keep the value as text instead of re-parsing it as a number,
# and resetting the Style afterwards avoids picking up an implicit "Text" format).

$cell = $ws.Range("D2")
$cell.Value = "'27.329.39"
$cell.Style = "Normal"
$cell = $ws.Range("E2")
$cell.Value = "  -0.77%  "
$cell.Style = "Normal"
$cell = $ws.Range("D3")
$cell.Value = "'1.710.74"
$cell.Style = "Normal"
$cell = $ws.Range("E3")
$cell.Value = "  -0.87%  "
$cell.Style = "Normal"
$cell = $ws.Range("E4")
$cell.Value = "  -0.01%  "
$cell.Style = "Normal"
$cell = $ws.Range("D5")
$cell.Value = "'224.72"
$cell.Style = "Normal"
$cell = $ws.Range("E5")
$cell.Value = "  -0.52%  "
$cell.Style = "Normal"
$cell = $ws.Range("E6")
$cell.Value = "  -1.08%  "
$cell.Style = "Normal"
$cell = $ws.Range("E7")
$cell.Value = "  +0.07%  "
$cell.Style = "Normal"
$cell = $ws.Range("D8")
$cell.Value = "'0.06679"
$cell.Style = "Normal"
$cell = $ws.Range("E8")
$cell.Value = "  +1.30%  "
$cell.Style = "Normal"
$cell = $ws.Range("E9")
$cell.Value = "  -0.46%  "
$cell.Style = "Normal"
$cell = $ws.Range("D10")
$cell.Value = "'20.81"
$cell.Style = "Normal"
$cell = $ws.Range("E10")
$cell.Value = "  -3.84%  "
$cell.Style = "Normal"
$cell = $ws.Range("D11")
$cell.Value = "'0.07684"
$cell.Style = "Normal"
$cell = $ws.Range("E11")
$cell.Value = "  -0.20%  "
$cell.Style = "Normal"
$cell = $ws.Range("D12")
$cell.Value = "'4.497"
$cell.Style = "Normal"
$cell = $ws.Range("E12")
$cell.Value = "  -2.53%  "
$cell.Style = "Normal"
$cell = $ws.Range("D13")
$cell.Value = "'1.946.28"
$cell.Style = "Normal"
$cell = $ws.Range("E13")
$cell.Value = "  -0.89%  "
$cell.Style = "Normal"
$cell = $ws.Range("D14")
$cell.Value = "'1.704.33"
$cell.Style = "Normal"
$cell = $ws.Range("E14")
$cell.Value = "  -1.15%  "
$cell.Style = "Normal"
$cell = $ws.Range("E15")
$cell.Value = "  +0.36%  "
$cell.Style = "Normal"
$cell = $ws.Range("D16")
$cell.Value = "'0.0₅8211"
$cell.Style = "Normal"
$cell = $ws.Range("E16")
$cell.Value = "  -0.94%  "
$cell.Style = "Normal"
$cell = $ws.Range("D17")
$cell.Value = "'67.90"
$cell.Style = "Normal"
$cell = $ws.Range("E17")
$cell.Value = "  -0.08%  "
$cell.Style = "Normal"
$cell = $ws.Range("D18")
$cell.Value = "'27.359.10"
$cell.Style = "Normal"
$cell = $ws.Range("E18")
$cell.Value = "  -0.67%  "
$cell.Style = "Normal"
$cell = $ws.Range("D19")
$cell.Value = "'221.97"
$cell.Style = "Normal"
$cell = $ws.Range("E19")
$cell.Value = "  +1.27%  "
$cell.Style = "Normal"
$cell = $ws.Range("E20")
$cell.Value = "  +0.01%  "
$cell.Style = "Normal"
$cell = $ws.Range("D21")
$cell.Value = "'4.643"
$cell.Style = "Normal"
$cell = $ws.Range("E21")
$cell.Value = "  -1.87%  "
$cell.Style = "Normal"
$cell = $ws.Range("D22")
$cell.Value = "'10.47"
$cell.Style = "Normal"
$cell = $ws.Range("E22")
$cell.Value = "  -1.40%  "
$cell.Style = "Normal"
$cell = $ws.Range("D23")
$cell.Value = "'6.017"
$cell.Style = "Normal"
$cell = $ws.Range("E23")
$cell.Value = "  -1.28%  "
$cell.Style = "Normal"
$cell = $ws.Range("E24")
$cell.Value = "  +0.06%  "
$cell.Style = "Normal"
$cell = $ws.Range("D25")
$cell.Value = "'144.95"
$cell.Style = "Normal"
$cell = $ws.Range("E25")
$cell.Value = "  -0.15%  "
$cell.Style = "Normal"
$cell = $ws.Range("E26")
$cell.Value = "  -4.94%  "
$cell.Style = "Normal"
$cell = $ws.Range("D27")
$cell.Value = "'0.1208"
$cell.Style = "Normal"
$cell = $ws.Range("E27")
$cell.Value = "  -2.24%  "
$cell.Style = "Normal"
$cell = $ws.Range("D28")
$cell.Value = "'7.235"
$cell.Style = "Normal"
$cell = $ws.Range("E28")
$cell.Value = "  -2.30%  "
$cell.Style = "Normal"
$cell = $ws.Range("D29")
$cell.Value = "'16.24"
$cell.Style = "Normal"
$cell = $ws.Range("E29")
$cell.Value = "  -1.97%  "
$cell.Style = "Normal"
$cell = $ws.Range("D30")
$cell.Value = "'0.05337"
$cell.Style = "Normal"
$cell = $ws.Range("E30")
$cell.Value = "  -3.28%  "
$cell.Style = "Normal"
$cell = $ws.Range("D31")
$cell.Value = "'1.293"
$cell.Style = "Normal"
$cell = $ws.Range("E31")
$cell.Value = "  -0.94%  "
$cell.Style = "Normal"
$cell = $ws.Range("D32")
$cell.Value = "'3.471"
$cell.Style = "Normal"
$cell = $ws.Range("E32")
$cell.Value = "  -2.72%  "
$cell.Style = "Normal"
$cell = $ws.Range("D33")
$cell.Value = "'3.424"
$cell.Style = "Normal"
$cell = $ws.Range("E33")
$cell.Value = "  -0.57%  "
$cell.Style = "Normal"
$cell = $ws.Range("D34")
$cell.Value = "'1.631"
$cell.Style = "Normal"
$cell = $ws.Range("E34")
$cell.Value = "  -1.72%  "
$cell.Style = "Normal"
$cell = $ws.Range("D35")
$cell.Value = "'2.871"
$cell.Style = "Normal"
$cell = $ws.Range("E35")
$cell.Value = "  +0.39%  "
$cell.Style = "Normal"
$cell = $ws.Range("D36")
$cell.Value = "'0.9537"
$cell.Style = "Normal"
$cell = $ws.Range("E36")
$cell.Value = "  -1.31%  "
$cell.Style = "Normal"
$cell = $ws.Range("E37")
$cell.Value = "  -1.35%  "
$cell.Style = "Normal"
$cell = $ws.Range("D38")
$cell.Value = "'0.5856"
$cell.Style = "Normal"
$cell = $ws.Range("E38")
$cell.Value = "  -2.07%  "
$cell.Style = "Normal"
$cell = $ws.Range("D39")
$cell.Value = "'1.148.22"
$cell.Style = "Normal"
$cell = $ws.Range("E39")
$cell.Value = "  +8.62%  "
$cell.Style = "Normal"
$cell = $ws.Range("E40")
$cell.Value = "  -1.03%  "
$cell.Style = "Normal"
$cell = $ws.Range("D41")
$cell.Value = "'5.790"
$cell.Style = "Normal"
$cell = $ws.Range("E41")
$cell.Value = "  -1.87%  "
$cell.Style = "Normal"
$cell = $ws.Range("D42")
$cell.Value = "'1.005"
$cell.Style = "Normal"
$cell = $ws.Range("E42")
$cell.Value = "  +0.11%  "
$cell.Style = "Normal"
$cell = $ws.Range("D43")
$cell.Value = "'0.8392"
$cell.Style = "Normal"
$cell = $ws.Range("E43")
$cell.Value = "  -1.85%  "
$cell.Style = "Normal"
$cell = $ws.Range("D44")
$cell.Value = "'100.97"
$cell.Style = "Normal"
$cell = $ws.Range("E44")
$cell.Value = "  -0.57%  "
$cell.Style = "Normal"
$cell = $ws.Range("D45")
$cell.Value = "'1.853.31"
$cell.Style = "Normal"
$cell = $ws.Range("E45")
$cell.Value = "  -0.86%  "
$cell.Style = "Normal"
$cell = $ws.Range("E46")
$cell.Value = "  -3.10%  "
$cell.Style = "Normal"
$cell = $ws.Range("D47")
$cell.Value = "'57.55"
$cell.Style = "Normal"
$cell = $ws.Range("E47")
$cell.Value = "  -2.36%  "
$cell.Style = "Normal"
$cell = $ws.Range("D48")
$cell.Value = "'0.4564"
$cell.Style = "Normal"
$cell = $ws.Range("E48")
$cell.Value = "  +2.12%  "
$cell.Style = "Normal"
$cell = $ws.Range("E49")
$cell.Value = "  -0.13%  "
$cell.Style = "Normal"
$cell = $ws.Range("D50")
$cell.Value = "'8.114"
$cell.Style = "Normal"
$cell = $ws.Range("E50")
$cell.Value = "  -1.19%  "
$cell.Style = "Normal"
$cell = $ws.Range("D51")
$cell.Value = "'0.05204"
$cell.Style = "Normal"
$cell = $ws.Range("E51")
$cell.Value = "  -0.74%  "
$cell.Style = "Normal"
